$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1251899
$ws.Range("I6").Value = 5000500
$ws.Range("J6").Value = 2365.3333
$ws.Range("K6").Value = 15001500
$ws.Range("L6").Value = 7095.999899999999
$ws.Range("M6").Value = -15001388
$ws.Range("N6").Value = -7319.999899999999
$ws.Range("H43").Value = 1217
$ws.Range("I43").Value = 1171.25
$ws.Range("J43").Value = 1400
$ws.Range("K43").Value = 1171.25
$ws.Range("L43").Value = 1400
$ws.Range("M43").Value = -1102.25
$ws.Range("N43").Value = -1538
$ws.Range("H116").Value = 1965.5883
$ws.Range("I116").Value = 1531.875
$ws.Range("J116").Value = 2351.111
$ws.Range("K116").Value = 1531.875
$ws.Range("L116").Value = 2351.111
$ws.Range("M116").Value = 1910.125
$ws.Range("N116").Value = -9235.111000000001
$ws.Range("H132").Value = 5323902
$ws.Range("I132").Value = 6415826
$ws.Range("J132").Value = 772.125
$ws.Range("K132").Value = 19247478
$ws.Range("L132").Value = 2316.375
$ws.Range("M132").Value = -19244948
$ws.Range("N132").Value = -7376.375
$ws.Range("H137").Value = 1372.9231
$ws.Range("I137").Value = 1174.1177
$ws.Range("J137").Value = 1469.4857
$ws.Range("K137").Value = 3522.3531
$ws.Range("L137").Value = 4408.4571
$ws.Range("M137").Value = -972.3531000000003
$ws.Range("N137").Value = -9508.4571
$ws.Range("H138").Value = 3757.8
$ws.Range("I138").Value = 1437.5483
$ws.Range("J138").Value = 5873.3237
$ws.Range("K138").Value = 4312.644899999999
$ws.Range("L138").Value = 17619.9711
$ws.Range("M138").Value = 827.3551000000007
$ws.Range("N138").Value = -27899.9711
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19898.955
$ws.Range("I32").Value = 3154.8904
$ws.Range("K32").Value = 3154.8904
$ws.Range("M32").Value = -2867.8904
$ws.Range("H45").Value = 1804.6666
$ws.Range("I45").Value = 1507.9333
$ws.Range("J45").Value = 2546.5
$ws.Range("K45").Value = 1507.9333
$ws.Range("L45").Value = 2546.5
$ws.Range("M45").Value = -1130.9333
$ws.Range("N45").Value = -3300.5
$ws.Range("H122").Value = 1973.421
$ws.Range("I122").Value = 1939.7222
$ws.Range("J122").Value = 2580
$ws.Range("K122").Value = 5819.1666
$ws.Range("L122").Value = 7740
$ws.Range("M122").Value = -3369.1666
$ws.Range("N122").Value = -12640
$ws.Range("H132").Value = 3851.2727
$ws.Range("I132").Value = 4203.2856
$ws.Range("J132").Value = 1880
$ws.Range("K132").Value = 12609.8568
$ws.Range("L132").Value = 5640
$ws.Range("M132").Value = -10079.8568
$ws.Range("N132").Value = -10700
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 22252266
$ws.Range("I107").Value = 25674604
$ws.Range("J107").Value = 7062
$ws.Range("K107").Value = 25674604
$ws.Range("L107").Value = 7062
$ws.Range("M107").Value = -25672684
$ws.Range("N107").Value = -10902
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1600
$ws.Range("I4").Value = 1500
$ws.Range("J4").Value = 1800
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = -1388
$ws.Range("N4").Value = -2024
$ws.Range("H31").Value = 16464.836
$ws.Range("I31").Value = 27985.838
$ws.Range("K31").Value = 27985.838
$ws.Range("M31").Value = -27690.838
$ws.Range("H34").Value = 16464.836
$ws.Range("I34").Value = 27985.838
$ws.Range("K34").Value = 27985.838
$ws.Range("M34").Value = -27783.838
$ws.Range("H58").Value = 8720.485000000001
$ws.Range("I58").Value = 1255.5652
$ws.Range("K58").Value = 1255.5652
$ws.Range("M58").Value = -1052.5652
$ws.Range("H122").Value = 2475.5925
$ws.Range("I122").Value = 2364.8635
$ws.Range("K122").Value = 7094.5905
$ws.Range("M122").Value = -4644.5905
$ws.Range("H136").Value = 8720.485000000001
$ws.Range("I136").Value = 1255.5652
$ws.Range("K136").Value = 3766.6956
$ws.Range("M136").Value = -1216.6956
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1895.8
$ws.Range("I41").Value = 793.3333
$ws.Range("J41").Value = 3549.5
$ws.Range("K41").Value = 2379.9999
$ws.Range("L41").Value = 10648.5
$ws.Range("M41").Value = -2041.9999
$ws.Range("N41").Value = -11324.5
$ws.Range("H68").Value = 1813.619
$ws.Range("J68").Value = 2292.2683
$ws.Range("L68").Value = 6876.804900000001
$ws.Range("N68").Value = -8498.804900000001
$ws.Range("H71").Value = 1813.619
$ws.Range("J71").Value = 2292.2683
$ws.Range("L71").Value = 20630.4147
$ws.Range("N71").Value = -28742.4147
$ws.Range("H86").Value = 1433.2858
$ws.Range("I86").Value = 1166.6666
$ws.Range("J86").Value = 1633.25
$ws.Range("K86").Value = 3499.9998
$ws.Range("L86").Value = 4899.75
$ws.Range("M86").Value = -2313.9998
$ws.Range("N86").Value = -7271.75
$ws.Range("H89").Value = 1433.2858
$ws.Range("I89").Value = 1166.6666
$ws.Range("J89").Value = 1633.25
$ws.Range("K89").Value = 10499.9994
$ws.Range("L89").Value = 14699.25
$ws.Range("M89").Value = -4571.999400000001
$ws.Range("N89").Value = -26555.25
$ws.Range("H107").Value = 240435.95
$ws.Range("I107").Value = 362.0645
$ws.Range("J107").Value = 516076.34
$ws.Range("K107").Value = 1086.1935
$ws.Range("L107").Value = 1548229.02
$ws.Range("M107").Value = 833.8064999999999
$ws.Range("N107").Value = -1552069.02
$ws.Range("H131").Value = 1599.5186
$ws.Range("J131").Value = 1603.171
$ws.Range("L131").Value = 4809.513
$ws.Range("N131").Value = -14889.513
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 7250.3335
$ws.Range("I41").Value = 4551
$ws.Range("J41").Value = 8600
$ws.Range("K41").Value = 4551
$ws.Range("L41").Value = 8600
$ws.Range("M41").Value = -4196
$ws.Range("N41").Value = -9310
$ws.Range("H70").Value = 47162.875
$ws.Range("I70").Value = 76773.86
$ws.Range("J70").Value = 5707.5
$ws.Range("K70").Value = 76773.86
$ws.Range("L70").Value = 5707.5
$ws.Range("M70").Value = -76503.86
$ws.Range("N70").Value = -6247.5
$ws.Range("H73").Value = 47162.875
$ws.Range("I73").Value = 76773.86
$ws.Range("J73").Value = 5707.5
$ws.Range("K73").Value = 76773.86
$ws.Range("L73").Value = 5707.5
$ws.Range("M73").Value = -75837.86
$ws.Range("N73").Value = -7579.5
$ws.Range("H102").Value = 242113.77
$ws.Range("I102").Value = 2424.8333
$ws.Range("J102").Value = 858456.7
$ws.Range("K102").Value = 2424.8333
$ws.Range("L102").Value = 858456.7
$ws.Range("M102").Value = -802.8332999999998
$ws.Range("N102").Value = -861700.7
$ws.Range("H132").Value = 2448.5151
$ws.Range("I132").Value = 2488
$ws.Range("J132").Value = 2357.7
$ws.Range("K132").Value = 7464
$ws.Range("L132").Value = 7073.099999999999
$ws.Range("M132").Value = -4934
$ws.Range("N132").Value = -12133.1
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1946.7742
$ws.Range("J2").Value = 1994.6428
$ws.Range("L2").Value = 1994.6428
$ws.Range("N2").Value = -2218.6428
$ws.Range("H7").Value = 1481.0769
$ws.Range("J7").Value = 2900
$ws.Range("L7").Value = 2900
$ws.Range("N7").Value = -3124
$ws.Range("H16").Value = 54063.156
$ws.Range("I16").Value = 84344.164
$ws.Range("J16").Value = 2152.8572
$ws.Range("K16").Value = 84344.164
$ws.Range("L16").Value = 2152.8572
$ws.Range("M16").Value = -84174.164
$ws.Range("N16").Value = -2492.8572
$ws.Range("H40").Value = 251500
$ws.Range("I40").Value = 501200
$ws.Range("J40").Value = 1800
$ws.Range("K40").Value = 501200
$ws.Range("L40").Value = 1800
$ws.Range("M40").Value = -501064
$ws.Range("N40").Value = -2072
$ws.Range("H61").Value = 3526.6667
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 4580
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 4580
$ws.Range("M61").Value = -2798
$ws.Range("N61").Value = -4984
$ws.Range("H93").Value = 1353.6875
$ws.Range("I93").Value = 1354.6666
$ws.Range("K93").Value = 1354.6666
$ws.Range("M93").Value = -106.6666
$ws.Range("H113").Value = 3526.6667
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 4580
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 4580
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -8920
$ws.Range("H126").Value = 1481.0769
$ws.Range("J126").Value = 2900
$ws.Range("L126").Value = 8700
$ws.Range("N126").Value = -13640
$ws.Range("H132").Value = 6269.4546
$ws.Range("I132").Value = 5896.4
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 17689.2
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -15159.2
$ws.Range("N132").Value = -35060
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 17977.889
$ws.Range("I2").Value = 14134
$ws.Range("J2").Value = 19899.834
$ws.Range("K2").Value = 14134
$ws.Range("L2").Value = 19899.834
$ws.Range("M2").Value = -14022
$ws.Range("N2").Value = -20123.834
$ws.Range("H132").Value = 3248.5
$ws.Range("I132").Value = 4061
$ws.Range("J132").Value = 2255.4443
$ws.Range("K132").Value = 12183
$ws.Range("L132").Value = 6766.3329
$ws.Range("M132").Value = -9653
$ws.Range("N132").Value = -11826.3329
$ws.Range("H136").Value = 710.83673
$ws.Range("J136").Value = 1249.7333
$ws.Range("L136").Value = 3749.199900000001
$ws.Range("N136").Value = -8849.1999
Write-Host "Applied all edits"
